$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "NC1-1AA - 007"

$ws.Range("A21").Value = 0.00000000047
$ws.Range("B21").Value = 1500
$ws.Range("C21").Value = 100
$ws.Range("D21").Value = 1257
$ws.Range("E21").Value = 433
$ws.Range("F21").Value = 1691

$ws.Range("A22").Value = 0.00000000022
$ws.Range("B22").Value = 1500
$ws.Range("C22").Value = 100
$ws.Range("D22").Value = 771
$ws.Range("E22").Value = 410
$ws.Range("F22").Value = 1186

$ws.Range("A21:A22").NumberFormat = "0.00E+00"

$ws.Range("D23").Select()
